# Mealy_R5 -> Mealy_R100 rename + refreshed result values (B/C/E at rows 14, 27, 40)
# and cursor moved onto the final result row (E40) on the renamed sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mealy_R5")

# Renaming the sheet also updates the workbook-level defined name
# ("Moore_R10_PDS") that points at it.
$ws.Name = "Mealy_R100"

# Updated measurement results (gcd-overflow fix changed the cached numbers).
$ws.Range("B14").Value = 3177
$ws.Range("C14").Value = 6028
$ws.Range("E14").Value = 24669

$ws.Range("B27").Value = 3177
$ws.Range("C27").Value = 6028
$ws.Range("E27").Value = 24669

$ws.Range("B40").Value = 3177
$ws.Range("C40").Value = 6028
$ws.Range("E40").Value = 6028

# Make the renamed sheet the active one, with the selection left on the
# last updated cell (matches the author re-saving after editing row 40).
$ws.Activate()
$ws.Range("E40").Select()
